# Apply "use case description edits" to the Use Case Description workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- 1. Update the "4.1." flow-of-activities text for "Registers in Application" ---
$ws.Range("C65").Value = "4.1. System will send an email authentication using the email input in basic data. (System activity 1.1.)"

# --- 2. Fill in the previously-empty 6th use case block (rows 70-78): "Insert Appointment Details" ---
$ws.Range("B70").Value = "Insert Appointment Details"
$ws.Range("B71").Value = "Dermatologist will record the events for the accomplished appointment"
$ws.Range("B72").Value = "Accomplished appointment"
$ws.Range("B73").Value = "The dermatologist will make a summary of all the necessary details and information on the accomplished appointment"
$ws.Range("B74").Value = "Dermatologist"
$ws.Range("B76").Value = "Dermatologist, Patient"
$ws.Range("B77").Value = "Accomplished Appointment"

# --- 3. Row height adjustments ---
$ws.Rows.Item(46).RowHeight = 30.75
$ws.Rows.Item(73).RowHeight = 30.75

# --- 4. View/selection changes ---
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("C10").Select()
